$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 782.3333
$ws.Range("I33").Value = 782.3333
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 782.3333
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -553.3333
$ws.Range("H137").Value = 2354.889
$ws.Range("I137").Value = 823.5
$ws.Range("J137").Value = 3580
$ws.Range("K137").Value = 2470.5
$ws.Range("L137").Value = 10740
$ws.Range("M137").Value = 79.5
$ws.Range("N137").Value = -15840
$ws.Range("H140").Value = 80780
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 80780
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 80780
$ws.Range("N140").Value = -91140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 882.44446
$ws.Range("I2").Value = 882.44446
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 882.44446
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -769.44446
$ws.Range("N2").Value = ""
$ws.Range("H32").Value = 963.7619
$ws.Range("I32").Value = 791.2778
$ws.Range("J32").Value = 1998.6666
$ws.Range("K32").Value = 791.2778
$ws.Range("L32").Value = 1998.6666
$ws.Range("M32").Value = -504.2778
$ws.Range("N32").Value = -2572.6666
$ws.Range("H45").Value = 4097.143
$ws.Range("I45").Value = 3314.3333
$ws.Range("J45").Value = 4684.25
$ws.Range("K45").Value = 3314.3333
$ws.Range("L45").Value = 4684.25
$ws.Range("M45").Value = -2937.3333
$ws.Range("N45").Value = -5438.25
$ws.Range("H74").Value = 14472.4
$ws.Range("I74").Value = 14472.4
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 14472.4
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -13598.4
$ws.Range("H77").Value = 14472.4
$ws.Range("I77").Value = 14472.4
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 72362
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -67994
$ws.Range("H116").Value = 882.44446
$ws.Range("I116").Value = 882.44446
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 882.44446
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1411.55554
$ws.Range("N116").Value = ""
$ws.Range("H132").Value = 4274.5454
$ws.Range("I132").Value = 5521.375
$ws.Range("J132").Value = 949.6667
$ws.Range("K132").Value = 16564.125
$ws.Range("L132").Value = 2849.0001
$ws.Range("M132").Value = -14034.125
$ws.Range("N132").Value = -7909.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 882.44446
$ws.Range("I3").Value = 882.44446
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 882.44446
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -768.44446
$ws.Range("N3").Value = ""
$ws.Range("H105").Value = 1696.4286
$ws.Range("I105").Value = 1581.8
$ws.Range("J105").Value = 1983
$ws.Range("K105").Value = 1581.8
$ws.Range("L105").Value = 1983
$ws.Range("M105").Value = 165.2
$ws.Range("N105").Value = -5477

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 701670.3
$ws.Range("I16").Value = 701670.3
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 701670.3
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -701383.3
$ws.Range("H31").Value = 8345.8
$ws.Range("I31").Value = 2470.6667
$ws.Range("J31").Value = 9814.583
$ws.Range("K31").Value = 2470.6667
$ws.Range("L31").Value = 9814.583
$ws.Range("M31").Value = -2175.6667
$ws.Range("N31").Value = -10404.583
$ws.Range("H34").Value = 8345.8
$ws.Range("I34").Value = 2470.6667
$ws.Range("J34").Value = 9814.583
$ws.Range("K34").Value = 2470.6667
$ws.Range("L34").Value = 9814.583
$ws.Range("M34").Value = -2268.6667
$ws.Range("N34").Value = -10218.583
$ws.Range("H58").Value = 2640.1538
$ws.Range("I58").Value = 1484.7273
$ws.Range("J58").Value = 8995
$ws.Range("K58").Value = 1484.7273
$ws.Range("L58").Value = 8995
$ws.Range("M58").Value = -1281.7273
$ws.Range("N58").Value = -9401
$ws.Range("H103").Value = 16491.5
$ws.Range("I103").Value = 16988.666
$ws.Range("J103").Value = 15000
$ws.Range("K103").Value = 16988.666
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = -15816.666
$ws.Range("N103").Value = -17344
$ws.Range("H113").Value = 701670.3
$ws.Range("I113").Value = 701670.3
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 701670.3
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -699500.3
$ws.Range("H136").Value = 2640.1538
$ws.Range("I136").Value = 1484.7273
$ws.Range("J136").Value = 8995
$ws.Range("K136").Value = 4454.1819
$ws.Range("L136").Value = 26985
$ws.Range("M136").Value = -1904.1819
$ws.Range("N136").Value = -32085
$ws.Range("H140").Value = 87890
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 87890
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 87890
$ws.Range("N140").Value = -98250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 197.5
$ws.Range("I61").Value = 200
$ws.Range("J61").Value = 195
$ws.Range("K61").Value = 600
$ws.Range("L61").Value = 585
$ws.Range("M61").Value = -385
$ws.Range("N61").Value = -1015
$ws.Range("H63").Value = 566.3333
$ws.Range("I63").Value = 566.3333
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1698.9999
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -949.9999
$ws.Range("H66").Value = 566.3333
$ws.Range("I66").Value = 566.3333
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 5096.9997
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -1352.9997
$ws.Range("H75").Value = 5999.3335
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 5999.3335
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 17998.0005
$ws.Range("N75").Value = -19994.0005
$ws.Range("H78").Value = 5999.3335
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 5999.3335
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 53994.0015
$ws.Range("N78").Value = -63978.0015
$ws.Range("H112").Value = 5000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 5000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 15000
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -17216
$ws.Range("H114").Value = 596.25
$ws.Range("I114").Value = 252.85715
$ws.Range("J114").Value = 3000
$ws.Range("K114").Value = 758.5714499999999
$ws.Range("L114").Value = 9000
$ws.Range("M114").Value = 2495.42855
$ws.Range("N114").Value = -15508
$ws.Range("H122").Value = 718.125
$ws.Range("I122").Value = 746
$ws.Range("J122").Value = 701.4
$ws.Range("K122").Value = 6714
$ws.Range("L122").Value = 6312.599999999999
$ws.Range("M122").Value = -4264
$ws.Range("N122").Value = -11212.6
$ws.Range("H131").Value = 2285
$ws.Range("I131").Value = 1950
$ws.Range("J131").Value = 2352
$ws.Range("K131").Value = 5850
$ws.Range("L131").Value = 7056
$ws.Range("M131").Value = -810
$ws.Range("N131").Value = -17136
$ws.Range("H137").Value = 5170
$ws.Range("I137").Value = 4975
$ws.Range("J137").Value = 5235
$ws.Range("K137").Value = 14925
$ws.Range("L137").Value = 15705
$ws.Range("M137").Value = -9825
$ws.Range("N137").Value = -25905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -12340
$ws.Range("H140").Value = 143728.67
$ws.Range("I140").Value = 165203
$ws.Range("J140").Value = 100780
$ws.Range("K140").Value = 165203
$ws.Range("L140").Value = 100780
$ws.Range("M140").Value = -160023
$ws.Range("N140").Value = -111140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2859.8948
$ws.Range("I61").Value = 775.7
$ws.Range("J61").Value = 5175.6665
$ws.Range("K61").Value = 775.7
$ws.Range("L61").Value = 5175.6665
$ws.Range("M61").Value = -573.7
$ws.Range("N61").Value = -5579.6665
$ws.Range("H76").Value = 22500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 22500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 22500
$ws.Range("N76").Value = -23176
$ws.Range("H79").Value = 22500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 22500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 22500
$ws.Range("N79").Value = -24840
$ws.Range("H93").Value = 1447.3334
$ws.Range("I93").Value = 1765.8
$ws.Range("J93").Value = 1049.25
$ws.Range("K93").Value = 1765.8
$ws.Range("L93").Value = 1049.25
$ws.Range("M93").Value = -517.8
$ws.Range("N93").Value = -3545.25
$ws.Range("H113").Value = 2859.8948
$ws.Range("I113").Value = 775.7
$ws.Range("J113").Value = 5175.6665
$ws.Range("K113").Value = 775.7
$ws.Range("L113").Value = 5175.6665
$ws.Range("M113").Value = 1394.3
$ws.Range("N113").Value = -9515.6665
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = ""
